# Apply weekly re-shuffle of Fruta/Hortaliza data rows (Jengibre subset).
# The rows keep their "static" columns (A,B,C,E,F,G,H,I,N,O,Q,R) but the
# "weekly" columns (D = Fecha, J = Volumen, K = Precio minimo,
# L = Precio maximo, M = Precio promedio ponderado, P = Precio $/Kg)
# get redistributed across rows 2..12 according to the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row -> source row (values to copy from the original layout)
$rowMap = @{
    2  = 7
    3  = 8
    4  = 12
    5  = 9
    6  = 11
    7  = 5
    8  = 6
    9  = 2
    10 = 10
    11 = 4
    12 = 3
}

# Snapshot the original values of the weekly columns before overwriting anything.
$orig = @{}
foreach ($r in 2..12) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

foreach ($r in 2..12) {
    $src = $orig[$rowMap[$r]]
    $ws.Cells.Item($r, 4).Value2  = $src.D
    $ws.Cells.Item($r, 10).Value2 = $src.J
    $ws.Cells.Item($r, 11).Value2 = $src.K
    $ws.Cells.Item($r, 12).Value2 = $src.L
    $ws.Cells.Item($r, 13).Value2 = $src.M
    $ws.Cells.Item($r, 16).Value2 = $src.P
}
